# collection_specialization_AE.xlsx edit:
# AE changed to normalized/denormalized layout; E31/E33 collection_group_id added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-18: "Vertical" free text format -> "Denormalized"
$rowsDenorm = 3..18
foreach ($r in $rowsDenorm) {
    $ws.Range("H$r").Value = "AE_DENORMALIZED"
    $ws.Range("I$r").Value = "Denormalized"
    $ws.Range("K$r").Value = "Adverse Event Free Text Format (Denormalized)"
}

# Rows 19-34: "Horizontal" free text format -> "Normalized"
$rowsNorm = 19..34
foreach ($r in $rowsNorm) {
    $ws.Range("H$r").Value = "AE_NORMALIZED"
    $ws.Range("I$r").Value = "Normalized"
    $ws.Range("K$r").Value = "Adverse Event Free Text Format (Normalized)"
}

# Add missing collection_group_id ("2-1") values on rows 31 and 33
$ws.Range("E31").Value = "2-1"
$ws.Range("E33").Value = "2-1"

# Minor font formatting nuance seen in original edit (cells H4:H18 picked up
# an explicit "applied font" style in the saved workbook).
$ws.Range("H4:H18").Font.Name = "Calibri"

Write-Output "done"
